$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 233 updates
# ---------------------------------------------------------------------------
$ws.Range("F233").Value = "BMC Software, Inc."
$ws.Range("G233").Value = "Compagnie intéressante"
$ws.Range("H233").Value = 2.5
$ws.Range("K233").Value = 45276.547962963
$ws.Range("O233").Value = "v1.1.0"

# ---------------------------------------------------------------------------
# Row 234 updates
# ---------------------------------------------------------------------------
$ws.Range("H234").Value = 3.5
$ws.Range("K234").Value = 45276.5418518518
$ws.Range("O234").Value = "v1.1.0"

# ---------------------------------------------------------------------------
# New rows 235-241: copy the formatting (and, for now, the values) of row 234
# -- our best template row for the new submissions -- into each new row, then
# overwrite with the real values. Columns that must not exist on a given row
# (per source data) are cleared afterwards so they don't leave an
# empty-but-styled cell behind.
# ---------------------------------------------------------------------------
for ($r = 235; $r -le 241; $r++) {
    $ws.Range("A234:Q234").Copy($ws.Range("A$r`:Q$r"))
}

# Columns that must stay completely empty on every new row (the template row
# carries column-level formatting for these that we don't want to keep).
$ws.Range("I235:I241").Clear()
$ws.Range("M235:M241").Clear()
$ws.Range("P235:P241").Clear()

# Row 235
$ws.Range("A235").Value = 235
$ws.Range("B235").Value = 1
$ws.Range("C235").Value = "GC"
$ws.Range("D235").Value = 45276
$ws.Range("E235").Value = 369
$ws.Range("F235").Value = "CDW Corporation"
$ws.Range("G235").Value = "CDW Corporation"
$ws.Range("H235").Value = 1
$ws.Range("J235").Value = $true
$ws.Range("K235").Value = 45276.511099537
$ws.Range("L235").Value = $false
$ws.Range("N235").Value = $false
$ws.Range("O235").Value = "v1.1.0"
$ws.Range("Q234").Copy($ws.Range("Q235"))

# Row 236
$ws.Range("A236").Value = 236
$ws.Range("B236").Value = 1
$ws.Range("C236").Value = "GC"
$ws.Range("D236").Value = 45276
$ws.Range("E236").Value = 1757
$ws.Range("F236").Value = "Telefonica Moviles S.A."
$ws.Range("G236").Value = "Telefonica"
$ws.Range("H236").Value = 2
$ws.Range("I234").Copy($ws.Range("I236"))
$ws.Range("I236").Value = "Telefonica"
$ws.Range("J236").Value = $true
$ws.Range("K236").Value = 45276.5113657407
$ws.Range("L236").Value = $false
$ws.Range("N236").Value = $false
$ws.Range("O236").Value = "v1.1.0"
$ws.Range("Q234").Copy($ws.Range("Q236"))

# Row 237
$ws.Range("A237").Value = 237
$ws.Range("B237").Value = 2
$ws.Range("C237").Value = "VG"
$ws.Range("D237").Value = 45276
$ws.Range("E237").Value = 570
$ws.Range("F237").Value = "Digital Insight Corporation"
$ws.Range("G237").Value = "Test du 16 décembre"
$ws.Range("H237").Value = 1
$ws.Range("J237").Value = $true
$ws.Range("K237").Value = 45276.5118287037
$ws.Range("L237").Value = $false
$ws.Range("N237").Value = $false
$ws.Range("O237").Value = "v1.1.0"
$ws.Range("Q234").Copy($ws.Range("Q237"))

# Row 238
$ws.Range("A238").Value = 238
$ws.Range("B238").Value = 2
$ws.Range("C238").Value = "VG"
$ws.Range("D238").Value = 45276
$ws.Range("E238").Value = 1263
$ws.Range("F238").Value = "Networks Associates, Inc."
$ws.Range("G238").Value = "Test du 16 décembre"
$ws.Range("H238").Value = 0.25
$ws.Range("J238").Value = $true
$ws.Range("K238").Value = 45276.512037037
$ws.Range("L238").Value = $false
$ws.Range("N238").Value = $false
$ws.Range("O238").Value = "v1.1.0"
$ws.Range("Q234").Copy($ws.Range("Q238"))

# Row 239
$ws.Range("A239").Value = 239
$ws.Range("B239").Value = 2
$ws.Range("C239").Value = "VG"
$ws.Range("D239").Value = 45276
$ws.Range("E239").Value = 1232
$ws.Range("F239").Value = "Murphy Oil Corporation"
$ws.Range("G239").Value = "Test du 16 décembre"
$ws.Range("H239").Value = 0.75
$ws.Range("J239").Value = $false
$ws.Range("K239").Value = 45276.5126041667
$ws.Range("L239").Value = $false
$ws.Range("N239").Value = $false
$ws.Range("O239").Value = "v1.1.0"
$ws.Range("Q234").Copy($ws.Range("Q239"))

# Row 240 (no company contact / no G or I column)
$ws.Range("A240").Value = 240
$ws.Range("B240").Value = 3
$ws.Range("C240").Value = "MFP"
$ws.Range("D240").Value = 45276
$ws.Range("E240").Value = 95
$ws.Range("F240").Value = "Altria Group"
$ws.Range("G240").Clear()
$ws.Range("H240").Value = 3
$ws.Range("J240").Value = $true
$ws.Range("K240").Value = 45276.5128472222
$ws.Range("L240").Value = $false
$ws.Range("N240").Value = $false
$ws.Range("O240").Value = "v1.1.0"
$ws.Range("Q234").Copy($ws.Range("Q240"))

# Row 241 (no Q column)
$ws.Range("A241").Value = 241
$ws.Range("B241").Value = 4
$ws.Range("C241").Value = "RMV"
$ws.Range("D241").Value = 45276
$ws.Range("E241").Value = 1554
$ws.Range("F241").Value = "Royal Caribbean Cruises Ltd."
$ws.Range("G241").Value = "Un beau voyage"
$ws.Range("H241").Value = 2
$ws.Range("J241").Value = $true
$ws.Range("K241").Value = 45276.5477662037
$ws.Range("L241").Value = $false
$ws.Range("N241").Value = $false
$ws.Range("O241").Value = "v1.1.0"
$ws.Range("Q241").Clear()

$excel.CutCopyMode = 0
